# GLI_No_FVC_Percent_Predicted.xlsx — trim the sample data down to a single
# "Example 1" row and relabel the spirometry headers with units.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two extra sample subjects (rows 3 and 4); row 2 becomes the only
# remaining data row and the sheet's used range shrinks to A1:F2.
$ws.Rows("3:4").Delete()

# Row 1 headers: clarify the spirometry columns with units.
$ws.Range("E1").Value = "FEV1 (L)"
$ws.Range("F1").Value = "FVC (L)"

# Row 2: replace the numeric "Unique ID" (1) with a descriptive label.
$ws.Range("A2").Value = "Example 1"

# Match the author's final selection/cursor position.
$ws.Range("E1").Select() | Out-Null
